$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume-change (E) columns for each coin row.
$ws.Range("D2").Value = "29.142.88"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "1.824.41"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.48"
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6168"
$ws.Range("E6").Value = "  -1.84%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9999"
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07335"
$ws.Range("E8").Value = "  -2.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2895"
$ws.Range("E9").Value = "  -1.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.94"
$ws.Range("E10").Value = "  -1.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07674"
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("D12").Value = "1.822.38"
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.955"
$ws.Range("E13").Value = "  -1.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6606"
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "81.94"
$ws.Range("E15").Value = "  -1.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008899"
$ws.Range("E16").Value = "  -5.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.832"
$ws.Range("E17").Value = "  -2.70%  "
$ws.Range("D18").Value = "29.123.63"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("D19").Value = "2.064.07"
$ws.Range("E19").Value = "  -0.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "237.64"
$ws.Range("E20").Value = "  +6.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.41"
$ws.Range("E21").Value = "  -1.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9999"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.156"
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("E24").Value = "  +0.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.85"
$ws.Range("E25").Value = "  -1.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1416"
$ws.Range("E26").Value = "  +1.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.420"
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("E28").Value = "  -1.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.484"
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05551"
$ws.Range("E30").Value = "  -3.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.092"
$ws.Range("E31").Value = "  -0.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.095"
$ws.Range("E32").Value = "  -1.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.205"
$ws.Range("E33").Value = "  +0.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.131"
$ws.Range("E36").Value = "  -0.57%  "
$ws.Range("E37").Value = "  -2.04%  "
$ws.Range("E38").Value = "  +2.91%  "
$ws.Range("D39").Value = "1.212.60"
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01756"
$ws.Range("E40").Value = "  -1.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.315"
$ws.Range("E41").Value = "  -2.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9163"
$ws.Range("E42").Value = "  +3.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9993"
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.00000000129"
$ws.Range("E44").Value = "  +5.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.36"
$ws.Range("E45").Value = "  -0.70%  "
$ws.Range("D46").Value = "1.971.61"
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "64.62"
$ws.Range("E47").Value = "  -1.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5086"
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4004"
$ws.Range("E49").Value = "  -1.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.026"
$ws.Range("E50").Value = "  +0.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05757"
$ws.Range("E51").Value = "  -1.11%  "

# Rows 34/35: ImmutableX and LidoDAOToken swapped rank positions.
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7359"
$ws.Range("E34").Value = "  -0.79%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.820"
$ws.Range("E35").Value = "  -0.75%  "

